# Test script - delete row 18 on sheet "nuevos 2025 "
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("nuevos 2025 ")
$ws.Rows("18:18").Delete()
